# previsao_retorno.xlsx update — "atualizei os dados da bibi e da add"
#
# Refreshes a handful of "INATIVO - X.Y meses sem comprar" counters (column J /
# situacao) that naturally ticked forward, one padrao_compra reclassification
# (column G, row 69), and three clients' recency metrics (rows 65, 69, 111 —
# prob_minima/prob_maxima/total_compras_historico/regularidade/ultima_compra/
# proxima_compra) on the Resumo_por_Cliente sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Resumo_por_Cliente")

# --- situacao (column J) "meses sem comprar" counter bumps ---------------
$ws.Range("J16").Value2  = "INATIVO - 1.4 meses sem comprar"
$ws.Range("J22").Value2  = "INATIVO - 28.1 meses sem comprar"
$ws.Range("J31").Value2  = "INATIVO - 6.5 meses sem comprar"
$ws.Range("J41").Value2  = "INATIVO - 6.4 meses sem comprar"
$ws.Range("J42").Value2  = "INATIVO - 24.7 meses sem comprar"
$ws.Range("J79").Value2  = "INATIVO - 21.2 meses sem comprar"
$ws.Range("J83").Value2  = "INATIVO - 20.7 meses sem comprar"
$ws.Range("J107").Value2 = "INATIVO - 20.3 meses sem comprar"

# --- padrao_compra (column G) reclassification ----------------------------
$ws.Range("G69").Value2 = "1x por mês - irregular (preferencialmente na 1ª quinzena)"

# --- row 65 (id_cliente 9247) refreshed metrics ----------------------------
$ws.Range("B65").Value2 = 0.25
$ws.Range("C65").Value2 = 0.17
$ws.Range("E65").Value2 = 29
$ws.Range("H65").Value2 = 45803.91217592593
$ws.Range("I65").Value2 = 45834.91217592593

# --- row 69 (id_cliente 9807) refreshed metrics ----------------------------
$ws.Range("C69").Value2 = 0.5
$ws.Range("D69").Value2 = 0.5
$ws.Range("E69").Value2 = 12
$ws.Range("F69").Value2 = 0.5
$ws.Range("H69").Value2 = 45803.93900462963
$ws.Range("I69").Value2 = 45834.93900462963

# --- row 111 (id_cliente 28458) refreshed metrics --------------------------
$ws.Range("E111").Value2 = 14861
$ws.Range("H111").Value2 = 45803.72020833333
$ws.Range("I111").Value2 = 45804.72020833333
